$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D are written as text,
# matching the original inline-string cell type (avoids Excel auto-converting
# values such as "3.017" or "314.03" into numbers).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2: Bitcoin
$ws.Range("D2").Value = '27.454.67'
$ws.Range("E2").Value = '  +2.31%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.878.43'
$ws.Range("E3").Value = '  +2.07%  '

# Row 4: TetherUSD
$ws.Range("D4").Value = '1.017'
$ws.Range("E4").Value = '  +0.91%  '

# Row 5: BNB
$ws.Range("D5").Value = '314.03'
$ws.Range("E5").Value = '  +1.60%  '

# Row 6: USDC
$ws.Range("D6").Value = '1.014'
$ws.Range("E6").Value = '  +0.73%  '

# Row 7: XRP
$ws.Range("D7").Value = '0.4796'
$ws.Range("E7").Value = '  +1.97%  '

# Row 8: Cardano
$ws.Range("D8").Value = '0.3772'
$ws.Range("E8").Value = '  +3.31%  '

# Row 9: Dogecoin
$ws.Range("D9").Value = '0.07434'
$ws.Range("E9").Value = '  +3.68%  '

# Row 10: Polygon
$ws.Range("D10").Value = '0.9425'
$ws.Range("E10").Value = '  +2.63%  '

# Row 11: Solana
$ws.Range("D11").Value = '20.90'
$ws.Range("E11").Value = '  +7.05%  '

# Row 12: TRON
$ws.Range("D12").Value = '0.07881'
$ws.Range("E12").Value = '  +3.73%  '

# Row 13: WrappedEther
$ws.Range("D13").Value = '1.869.77'
$ws.Range("E13").Value = '  +0.16%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '5.457'
$ws.Range("E14").Value = '  +3.38%  '

# Row 15: Chainlink
$ws.Range("D15").Value = '6.616'
$ws.Range("E15").Value = '  +3.31%  '

# Row 16: Litecoin
$ws.Range("D16").Value = '90.83'
$ws.Range("E16").Value = '  +3.55%  '

# Row 17: BinanceUSD
$ws.Range("D17").Value = '1.016'
$ws.Range("E17").Value = '  +0.60%  '

# Row 18: ShibaInu
$ws.Range("D18").Value = '0.000008900'
$ws.Range("E18").Value = '  +3.29%  '

# Row 19: Dai
$ws.Range("D19").Value = '1.014'
$ws.Range("E19").Value = '  +0.72%  '

# Row 20: Avalanche
$ws.Range("D20").Value = '14.94'
$ws.Range("E20").Value = '  +2.96%  '

# Row 21: WrappedBTC
$ws.Range("D21").Value = '27.477.61'
$ws.Range("E21").Value = '  +2.24%  '

# Row 22: Uniswap
$ws.Range("D22").Value = '5.152'
$ws.Range("E22").Value = '  +2.99%  '

# Row 23: Cosmos
$ws.Range("E23").Value = '  +1.24%  '

# Row 24: Toncoin
$ws.Range("D24").Value = '1.958'
$ws.Range("E24").Value = '  +1.62%  '

# Row 25: Monero
$ws.Range("D25").Value = '154.00'
$ws.Range("E25").Value = '  +1.54%  '

# Row 26: EthereumClassic
$ws.Range("D26").Value = '18.66'
$ws.Range("E26").Value = '  +3.02%  '

# Row 27: LidoDAOToken
$ws.Range("D27").Value = '2.037'
$ws.Range("E27").Value = '  +2.07%  '

# Row 28: BitcoinCash
$ws.Range("D28").Value = '116.24'
$ws.Range("E28").Value = '  +2.05%  '

# Row 29: InternetComputer(DFINITY)
$ws.Range("D29").Value = '5.043'
$ws.Range("E29").Value = '  +4.05%  '

# Row 30: Stellar
$ws.Range("D30").Value = '0.08941'
$ws.Range("E30").Value = '  +1.40%  '

# Row 31: HuobiToken
$ws.Range("D31").Value = '3.329'
$ws.Range("E31").Value = '  +0.65%  '

# Row 32: ARBITRUM
$ws.Range("D32").Value = '1.226'
$ws.Range("E32").Value = '  +4.92%  '

# Row 33: Filecoin
$ws.Range("D33").Value = '4.604'
$ws.Range("E33").Value = '  +2.82%  '

# Row 34: ImmutableX
$ws.Range("D34").Value = '0.7484'
$ws.Range("E34").Value = '  +1.03%  '

# Row 35: RenderToken
$ws.Range("D35").Value = '2.728'
$ws.Range("E35").Value = '  -0.46%  '

# Row 36: VeChain
$ws.Range("D36").Value = '0.02065'
$ws.Range("E36").Value = '  +6.30%  '

# Row 37: TrustWalletToken
$ws.Range("D37").Value = '1.126'
$ws.Range("E37").Value = '  +3.72%  '

# Row 38: Hedera
$ws.Range("D38").Value = '0.05311'
$ws.Range("E38").Value = '  +1.18%  '

# Row 39: MXToken
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '3.017'
$ws.Range("E39").Value = '  +1.35%  '

# Row 40: TheSandbox
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.5396'
$ws.Range("E40").Value = '  +4.03%  '

# Row 41: FraxShare
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '7.109'
$ws.Range("E41").Value = '  +3.35%  '

# Row 42: Algorand
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.1537'
$ws.Range("E42").Value = '  +1.83%  '

# Row 43: Aptos
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '8.486'
$ws.Range("E43").Value = '  +4.25%  '

# Row 44: EnergySwap
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '10.68'
$ws.Range("E44").Value = '  +1.72%  '

# Row 45: Decentraland
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '0.4856'
$ws.Range("E45").Value = '  +3.67%  '

# Row 46: PaxDollar
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '1.015'
$ws.Range("E46").Value = '  +0.71%  '

# Row 47: NEARProtocol
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '1.670'
$ws.Range("E47").Value = '  +4.76%  '

# Row 48: Quant
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '103.50'
$ws.Range("E48").Value = '  +1.92%  '

# Row 49: Aave
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '67.20'
$ws.Range("E49").Value = '  +2.68%  '

# Row 50: Cronos
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.06121'
$ws.Range("E50").Value = '  +1.46%  '

# Row 51: EOS
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").Value = '0.9021'
$ws.Range("E51").Value = '  +1.85%  '

# Restore default styling on column D (remove the temporary text format)
$ws.Range("D2:D51").Style = "Normal"
